$d = $word.ActiveDocument

# --- Table 1 (Java): add gradientBitmap row at the end ---
$t1 = $d.Tables(1)
$row1 = $t1.Rows.Add()
$idx1 = $row1.Index
$t1.Cell($idx1, 1).Range.Text = "gradientBitmap"
$t1.Cell($idx1, 2).Range.Text = "27.82"
$t1.Cell($idx1, 3).Range.Text = "55.13"

# --- Table 3 (Dart): add gradientBitmap row at the end ---
$t3 = $d.Tables(3)
$row3 = $t3.Rows.Add()
$idx3 = $row3.Index
$t3.Cell($idx3, 1).Range.Text = "gradientBitmap"
$t3.Cell($idx3, 2).Range.Text = "45.124"
$t3.Cell($idx3, 3).Range.Text = "104.816"
